$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns F and G, styled like the other headers (copy style from E1)
$ws.Range("F1").Value = "annuity_level"
$ws.Range("G1").Value = "premium_leveled"
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)  # xlPasteFormats

# Updated values for columns C, D, E and new columns F, G
$ws.Range("C2").Value = 13.59032634496769
$ws.Range("D2").Value = 13.35585960525418
$ws.Range("E2").Value = 3450.494620393525
$ws.Range("F2").Value = 13.92317044246741
$ws.Range("G2").Value = 247.8239158711355

$ws.Range("C3").Value = 13.59032634496769
$ws.Range("D3").Value = 13.44827497171428
$ws.Range("E3").Value = 2090.47773611435
$ws.Range("F3").Value = 14.00441067997329
$ws.Range("G3").Value = 149.2728101085891

$ws.Range("C4").Value = 13.59032634496769
$ws.Range("D4").Value = 13.32037861615614
$ws.Range("E4").Value = 3972.645276638446
$ws.Range("F4").Value = 13.8905593294462
$ws.Range("G4").Value = 285.9960626795603
